# Applies the Pandaemonium Profits scheduled-runner update.
# For each affected leve row, updates the live market-price / profit
# columns (H:N) to the newly refreshed values. Where a profit column
# no longer applies for a row, its cell is cleared instead of zeroed,
# matching the upstream data export behaviour.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 82
$ws.Range("H82").Value = 980
$ws.Range("I82").Value = 980
$ws.Range("K82").Value = 2940
$ws.Range("M82").Value = -2534
# Row 85
$ws.Range("H85").Value = 980
$ws.Range("I85").Value = 980
$ws.Range("K85").Value = 2940
$ws.Range("M85").Value = -1536
# Row 116
$ws.Range("H116").Value = 3393.4666
$ws.Range("I116").Value = 3233.3333
$ws.Range("J116").Value = 3633.6667
$ws.Range("K116").Value = 3233.3333
$ws.Range("L116").Value = 3633.6667
$ws.Range("M116").Value = 208.6667000000002
$ws.Range("N116").Value = -10517.6667

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 854.7931
$ws.Range("I2").Value = 995.64703
$ws.Range("J2").Value = 655.25
$ws.Range("K2").Value = 995.64703
$ws.Range("L2").Value = 655.25
$ws.Range("M2").Value = -882.64703
$ws.Range("N2").Value = -881.25
# Row 74
$ws.Range("H74").Value = 4482.613
$ws.Range("I74").Value = 5126.423
$ws.Range("J74").Value = 1134.8
$ws.Range("K74").Value = 5126.423
$ws.Range("L74").Value = 1134.8
$ws.Range("M74").Value = -4252.423
$ws.Range("N74").Value = -2882.8
# Row 77
$ws.Range("H77").Value = 4482.613
$ws.Range("I77").Value = 5126.423
$ws.Range("J77").Value = 1134.8
$ws.Range("K77").Value = 25632.115
$ws.Range("L77").Value = 5674
$ws.Range("M77").Value = -21264.115
$ws.Range("N77").Value = -14410
# Row 116
$ws.Range("H116").Value = 854.7931
$ws.Range("I116").Value = 995.64703
$ws.Range("J116").Value = 655.25
$ws.Range("K116").Value = 995.64703
$ws.Range("L116").Value = 655.25
$ws.Range("M116").Value = 1298.35297
$ws.Range("N116").Value = -5243.25
# Row 122
$ws.Range("H122").Value = 1860.2084
$ws.Range("I122").Value = 1691.6154
$ws.Range("J122").Value = 2059.4546
$ws.Range("K122").Value = 5074.8462
$ws.Range("L122").Value = 6178.3638
$ws.Range("M122").Value = -2624.8462
$ws.Range("N122").Value = -11078.3638

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 854.7931
$ws.Range("I3").Value = 995.64703
$ws.Range("J3").Value = 655.25
$ws.Range("K3").Value = 995.64703
$ws.Range("L3").Value = 655.25
$ws.Range("M3").Value = -881.64703
$ws.Range("N3").Value = -883.25
# Row 132
$ws.Range("H132").Value = 62922.082
$ws.Range("J132").Value = 62922.082
$ws.Range("L132").Value = 62922.082
$ws.Range("N132").Value = -73042.08199999999

$ws = $wb.Worksheets.Item("CRP")
# Row 3
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("N3").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 7589
$ws.Range("I3").Value = 4569.643
$ws.Range("J3").Value = 12872.875
$ws.Range("K3").Value = 13708.929
$ws.Range("L3").Value = 38618.625
$ws.Range("M3").Value = -13596.929
$ws.Range("N3").Value = -38842.625
# Row 7
$ws.Range("H7").Value = 1900
$ws.Range("I7").Value = 1900
$ws.Range("J7").Value = 1900
$ws.Range("K7").Value = 5700
$ws.Range("L7").Value = 5700
$ws.Range("M7").Value = -5588
$ws.Range("N7").Value = -5924
# Row 34
$ws.Range("H34").Value = 2747.56
$ws.Range("J34").Value = 3526.7896
$ws.Range("L34").Value = 10580.3688
$ws.Range("N34").Value = -10748.3688
# Row 92
$ws.Range("H92").Value = 769.4167
$ws.Range("I92").Value = 588.8333
$ws.Range("J92").Value = 950
$ws.Range("K92").Value = 1766.4999
$ws.Range("L92").Value = 2850
$ws.Range("M92").Value = -518.4999
$ws.Range("N92").Value = -5346
# Row 109
$ws.Range("H109").Value = 1964.55
$ws.Range("I109").Value = 623.875
$ws.Range("J109").Value = 2858.3333
$ws.Range("K109").Value = 1871.625
$ws.Range("L109").Value = 8574.999899999999
$ws.Range("M109").Value = -831.625
$ws.Range("N109").Value = -10654.9999
# Row 112
$ws.Range("H112").Value = 1086.8334
$ws.Range("I112").Value = 1086.8334
$ws.Range("J112").Value = 0
$ws.Range("K112").Value = 3260.5002
$ws.Range("L112").Value = 0
$ws.Range("M112").Value = -2152.5002
$ws.Range("N112").ClearContents()
# Row 121
$ws.Range("H121").Value = 950
$ws.Range("I121").Value = 400
$ws.Range("J121").Value = 1500
$ws.Range("K121").Value = 1200
$ws.Range("L121").Value = 4500
$ws.Range("M121").Value = 110
$ws.Range("N121").Value = -7120

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 958.3333
$ws.Range("I22").Value = 1100
$ws.Range("J22").Value = 887.5
$ws.Range("K22").Value = 1100
$ws.Range("L22").Value = 887.5
$ws.Range("M22").Value = -805
$ws.Range("N22").Value = -1477.5
# Row 27
$ws.Range("H27").Value = 958.3333
$ws.Range("I27").Value = 1100
$ws.Range("J27").Value = 887.5
$ws.Range("K27").Value = 1100
$ws.Range("L27").Value = 887.5
$ws.Range("M27").Value = -993
$ws.Range("N27").Value = -1101.5
# Row 61
$ws.Range("H61").Value = 11683.345
$ws.Range("I61").Value = 15046.223
$ws.Range("K61").Value = 15046.223
$ws.Range("M61").Value = -14844.223
# Row 70
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
# Row 73
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
# Row 113
$ws.Range("H113").Value = 11683.345
$ws.Range("I113").Value = 15046.223
$ws.Range("K113").Value = 15046.223
$ws.Range("M113").Value = -12876.223
# Row 122
$ws.Range("H122").Value = 5928.6333
$ws.Range("I122").Value = 5534.36
$ws.Range("J122").Value = 7900
$ws.Range("K122").Value = 16603.08
$ws.Range("L122").Value = 23700
$ws.Range("M122").Value = -14153.08
$ws.Range("N122").Value = -28600
# Row 135
$ws.Range("H135").Value = 37604.285
$ws.Range("J135").Value = 37604.285
$ws.Range("L135").Value = 37604.285
$ws.Range("N135").Value = -47744.285
# Row 136
$ws.Range("H136").Value = 5001.6665
$ws.Range("I136").Value = 3925.1853
$ws.Range("J136").Value = 6616.3887
$ws.Range("K136").Value = 11775.5559
$ws.Range("L136").Value = 19849.1661
$ws.Range("M136").Value = -9225.555899999999
$ws.Range("N136").Value = -24949.1661

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 1813.1154
$ws.Range("I122").Value = 1196.8636
$ws.Range("J122").Value = 5202.5
$ws.Range("K122").Value = 3590.5908
$ws.Range("L122").Value = 15607.5
$ws.Range("M122").Value = -1140.5908
$ws.Range("N122").Value = -20507.5
